$d = $word.ActiveDocument
$apos = [char]0x2019

# 1) Merge the runs of "C'est ok de dire la meme explication que partie d'avant?"
#    into a single run (also drops the proofErr spell/gram markers in between).
$find1 = " C" + $apos + "est ok de dire la meme explication que partie d" + $apos + "avant?"
$d.Content.Find.Execute($find1, $false, $false, $false, $false, $false, $true, 1, $false, $find1, 2) | Out-Null

# 2) Merge " Part " + bookmark(_GoBack) + "2" into a single run " Part 2"
#    (this also removes the _GoBack bookmark that previously sat here).
$find2 = " Part 2"
$d.Content.Find.Execute($find2, $false, $false, $false, $false, $false, $true, 1, $false, $find2, 2) | Out-Null

# 3) Merge the runs of "Same partition si N ou labels? Per class?????"
$find3 = "Same partition si N ou labels? Per class?????"
$d.Content.Find.Execute($find3, $false, $false, $false, $false, $false, $true, 1, $false, $find3, 2) | Out-Null

# 4) Merge the runs of "Comment faire mean et std?"
$find4 = "Comment faire mean et std?"
$d.Content.Find.Execute($find4, $false, $false, $false, $false, $false, $true, 1, $false, $find4, 2) | Out-Null

# 5) Merge the runs of "Faire un ttest pour partition?"
$find5 = "Faire un ttest pour partition?"
$d.Content.Find.Execute($find5, $false, $false, $false, $false, $false, $true, 1, $false, $find5, 2) | Out-Null

# 6) Fill the trailing empty list paragraph with a new question, and append a
#    brand-new list paragraph (same numbering) with a second question that has
#    a "_GoBack" bookmark sitting between its two halves.
$lastIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastIndex)
$r = $p.Range
$r.InsertAfter("How are the confusion matrix and the classification error linked? How do they differ?")

$p = $d.Paragraphs.Item($lastIndex)
$r = $p.Range
$r.Font.NameAscii = "Cambria"
$r.Font.Name = "Cambria"
$r.Font.NameBi = "Times New Roman"
$r.InsertParagraphAfter()

$newIndex = $d.Paragraphs.Count
$p2 = $d.Paragraphs.Item($newIndex)
$r2 = $p2.Range
$r2.InsertAfter("Can you find a way of modifying your classifier to give more weight to a certain class?")

$p2 = $d.Paragraphs.Item($newIndex)
$r2 = $p2.Range
$r2.Font.NameAscii = "Cambria"
$r2.Font.Name = "Cambria"
$r2.Font.NameBi = "Times New Roman"

$splitRange = $r2.Duplicate
$splitRange.Find.Execute("classif") | Out-Null
$splitRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null

Write-Output "done"
